$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(112, 8).Value = 5343.3335  # H112
$ws.Cells.Item(112, 10).Value = 5961.25  # J112
$ws.Cells.Item(112, 12).Value = 17883.75  # L112
$ws.Cells.Item(112, 14).Value = -20099.75  # N112
$ws.Cells.Item(113, 8).Value = 15029.889  # H113
$ws.Cells.Item(113, 9).Value = 16471.125  # I113
$ws.Cells.Item(113, 11).Value = 16471.125  # K113
$ws.Cells.Item(113, 13).Value = -13217.125  # M113
$ws.Cells.Item(137, 8).Value = 1483.6364  # H137
$ws.Cells.Item(137, 9).Value = 983.9167  # I137
$ws.Cells.Item(137, 10).Value = 2083.3  # J137
$ws.Cells.Item(137, 11).Value = 2951.7501  # K137
$ws.Cells.Item(137, 12).Value = 6249.900000000001  # L137
$ws.Cells.Item(137, 13).Value = -401.7501000000002  # M137
$ws.Cells.Item(137, 14).Value = -11349.9  # N137

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 0  # H5
$ws.Cells.Item(5, 9).Value = 0  # I5
$ws.Cells.Item(5, 11).Value = 0  # K5
$ws.Cells.Item(5, 13).ClearContents()  # M5
$ws.Cells.Item(32, 8).Value = 3812.3704  # H32
$ws.Cells.Item(32, 9).Value = 2461.0952  # I32
$ws.Cells.Item(32, 11).Value = 2461.0952  # K32
$ws.Cells.Item(32, 13).Value = -2174.0952  # M32
$ws.Cells.Item(61, 8).Value = 2101.6453  # H61
$ws.Cells.Item(61, 9).Value = 1350.2142  # I61
$ws.Cells.Item(61, 11).Value = 1350.2142  # K61
$ws.Cells.Item(61, 13).Value = -1138.2142  # M61
$ws.Cells.Item(64, 8).Value = 0  # H64
$ws.Cells.Item(64, 9).Value = 0  # I64
$ws.Cells.Item(64, 11).Value = 0  # K64
$ws.Cells.Item(64, 13).ClearContents()  # M64
$ws.Cells.Item(67, 8).Value = 0  # H67
$ws.Cells.Item(67, 9).Value = 0  # I67
$ws.Cells.Item(67, 11).Value = 0  # K67
$ws.Cells.Item(67, 13).ClearContents()  # M67
$ws.Cells.Item(74, 8).Value = 429.33334  # H74
$ws.Cells.Item(74, 9).Value = 429.33334  # I74
$ws.Cells.Item(74, 11).Value = 429.33334  # K74
$ws.Cells.Item(74, 13).Value = 444.66666  # M74
$ws.Cells.Item(77, 8).Value = 429.33334  # H77
$ws.Cells.Item(77, 9).Value = 429.33334  # I77
$ws.Cells.Item(77, 11).Value = 2146.6667  # K77
$ws.Cells.Item(77, 13).Value = 2221.3333  # M77
$ws.Cells.Item(122, 8).Value = 1641.7587  # H122
$ws.Cells.Item(122, 9).Value = 1649.619  # I122
$ws.Cells.Item(122, 11).Value = 4948.857  # K122
$ws.Cells.Item(122, 13).Value = -2498.857  # M122
$ws.Cells.Item(136, 8).Value = 2101.6453  # H136
$ws.Cells.Item(136, 9).Value = 1350.2142  # I136
$ws.Cells.Item(136, 11).Value = 4050.6426  # K136
$ws.Cells.Item(136, 13).Value = -1500.6426  # M136

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 0  # H4
$ws.Cells.Item(4, 9).Value = 0  # I4
$ws.Cells.Item(4, 11).Value = 0  # K4
$ws.Cells.Item(4, 13).ClearContents()  # M4
$ws.Cells.Item(107, 8).Value = 2592.5  # H107
$ws.Cells.Item(107, 9).Value = 2135.7144  # I107
$ws.Cells.Item(107, 11).Value = 2135.7144  # K107
$ws.Cells.Item(107, 13).Value = -215.7143999999998  # M107
$ws.Cells.Item(108, 8).Value = 34997.5  # H108
$ws.Cells.Item(108, 10).Value = 34997.5  # J108
$ws.Cells.Item(108, 12).Value = 34997.5  # L108
$ws.Cells.Item(108, 14).Value = -42677.5  # N108
$ws.Cells.Item(122, 8).Value = 0  # H122
$ws.Cells.Item(122, 10).Value = 0  # J122
$ws.Cells.Item(122, 12).Value = 0  # L122
$ws.Cells.Item(122, 14).ClearContents()  # N122
$ws.Cells.Item(134, 8).Value = 7134.6387  # H134
$ws.Cells.Item(134, 9).Value = 8638.521000000001  # I134
$ws.Cells.Item(134, 10).Value = 4473.923  # J134
$ws.Cells.Item(134, 11).Value = 25915.563  # K134
$ws.Cells.Item(134, 12).Value = 13421.769  # L134
$ws.Cells.Item(134, 13).Value = -23380.563  # M134
$ws.Cells.Item(134, 14).Value = -18491.769  # N134

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3158.5908  # H31
$ws.Cells.Item(31, 9).Value = 1437.7693  # I31
$ws.Cells.Item(31, 10).Value = 5644.222  # J31
$ws.Cells.Item(31, 11).Value = 1437.7693  # K31
$ws.Cells.Item(31, 12).Value = 5644.222  # L31
$ws.Cells.Item(31, 13).Value = -1142.7693  # M31
$ws.Cells.Item(31, 14).Value = -6234.222  # N31
$ws.Cells.Item(34, 8).Value = 3158.5908  # H34
$ws.Cells.Item(34, 9).Value = 1437.7693  # I34
$ws.Cells.Item(34, 10).Value = 5644.222  # J34
$ws.Cells.Item(34, 11).Value = 1437.7693  # K34
$ws.Cells.Item(34, 12).Value = 5644.222  # L34
$ws.Cells.Item(34, 13).Value = -1235.7693  # M34
$ws.Cells.Item(34, 14).Value = -6048.222  # N34
$ws.Cells.Item(58, 8).Value = 1141.0769  # H58
$ws.Cells.Item(58, 9).Value = 798.8125  # I58
$ws.Cells.Item(58, 10).Value = 1688.7  # J58
$ws.Cells.Item(58, 11).Value = 798.8125  # K58
$ws.Cells.Item(58, 12).Value = 1688.7  # L58
$ws.Cells.Item(58, 13).Value = -595.8125  # M58
$ws.Cells.Item(58, 14).Value = -2094.7  # N58
$ws.Cells.Item(81, 8).Value = 27750  # H81
$ws.Cells.Item(81, 10).Value = 27750  # J81
$ws.Cells.Item(81, 12).Value = 27750  # L81
$ws.Cells.Item(81, 14).Value = -29746  # N81
$ws.Cells.Item(84, 8).Value = 27750  # H84
$ws.Cells.Item(84, 10).Value = 27750  # J84
$ws.Cells.Item(84, 12).Value = 83250  # L84
$ws.Cells.Item(84, 14).Value = -93234  # N84
$ws.Cells.Item(136, 8).Value = 1141.0769  # H136
$ws.Cells.Item(136, 9).Value = 798.8125  # I136
$ws.Cells.Item(136, 10).Value = 1688.7  # J136
$ws.Cells.Item(136, 11).Value = 2396.4375  # K136
$ws.Cells.Item(136, 12).Value = 5066.1  # L136
$ws.Cells.Item(136, 13).Value = 153.5625  # M136
$ws.Cells.Item(136, 14).Value = -10166.1  # N136
$ws.Cells.Item(141, 8).Value = 24396.6  # H141
$ws.Cells.Item(141, 10).Value = 24396.6  # J141
$ws.Cells.Item(141, 12).Value = 24396.6  # L141
$ws.Cells.Item(141, 14).Value = -34756.6  # N141

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(103, 8).Value = 1360.25  # H103
$ws.Cells.Item(103, 9).Value = 1350  # I103
$ws.Cells.Item(103, 10).Value = 1370.5  # J103
$ws.Cells.Item(103, 11).Value = 4050  # K103
$ws.Cells.Item(103, 12).Value = 4111.5  # L103
$ws.Cells.Item(103, 13).Value = -3171  # M103
$ws.Cells.Item(103, 14).Value = -5869.5  # N103
$ws.Cells.Item(107, 8).Value = 530.4286  # H107
$ws.Cells.Item(107, 9).Value = 260  # I107
$ws.Cells.Item(107, 11).Value = 780  # K107
$ws.Cells.Item(107, 13).Value = 1140  # M107
$ws.Cells.Item(131, 8).Value = 5216387  # H131
$ws.Cells.Item(131, 10).Value = 8379.098  # J131
$ws.Cells.Item(131, 12).Value = 25137.294  # L131
$ws.Cells.Item(131, 14).Value = -35217.294  # N131
$ws.Cells.Item(134, 8).Value = 1716.2174  # H134
$ws.Cells.Item(134, 9).Value = 1374  # I134
$ws.Cells.Item(134, 10).Value = 3997.6667  # J134
$ws.Cells.Item(134, 11).Value = 4122  # K134
$ws.Cells.Item(134, 12).Value = 11993.0001  # L134
$ws.Cells.Item(134, 13).Value = 948  # M134
$ws.Cells.Item(134, 14).Value = -22133.0001  # N134
$ws.Cells.Item(137, 8).Value = 5357.154  # H137
$ws.Cells.Item(137, 9).Value = 3528.3333  # I137
$ws.Cells.Item(137, 10).Value = 6924.7144  # J137
$ws.Cells.Item(137, 11).Value = 10584.9999  # K137
$ws.Cells.Item(137, 12).Value = 20774.1432  # L137
$ws.Cells.Item(137, 13).Value = -5484.999899999999  # M137
$ws.Cells.Item(137, 14).Value = -30974.1432  # N137

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2113.5518  # H102
$ws.Cells.Item(102, 9).Value = 2061.0557  # I102
$ws.Cells.Item(102, 11).Value = 2061.0557  # K102
$ws.Cells.Item(102, 13).Value = -439.0556999999999  # M102
$ws.Cells.Item(126, 8).Value = 58369.5  # H126
$ws.Cells.Item(126, 10).Value = 127530  # J126
$ws.Cells.Item(126, 12).Value = 382590  # L126
$ws.Cells.Item(126, 14).Value = -387530  # N126
$ws.Cells.Item(133, 8).Value = 25000  # H133
$ws.Cells.Item(133, 10).Value = 25000  # J133
$ws.Cells.Item(133, 12).Value = 25000  # L133
$ws.Cells.Item(133, 14).Value = -35120  # N133
$ws.Cells.Item(138, 8).Value = 22611.9  # H138
$ws.Cells.Item(138, 10).Value = 25000  # J138
$ws.Cells.Item(138, 12).Value = 25000  # L138
$ws.Cells.Item(138, 14).Value = -35280  # N138

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 1662.8889  # H46
$ws.Cells.Item(46, 9).Value = 1233  # I46
$ws.Cells.Item(46, 11).Value = 1233  # K46
$ws.Cells.Item(46, 13).Value = -1045  # M46
$ws.Cells.Item(61, 8).Value = 3562.375  # H61
$ws.Cells.Item(61, 9).Value = 3083.1667  # I61
$ws.Cells.Item(61, 11).Value = 3083.1667  # K61
$ws.Cells.Item(61, 13).Value = -2881.1667  # M61
$ws.Cells.Item(113, 8).Value = 3562.375  # H113
$ws.Cells.Item(113, 9).Value = 3083.1667  # I113
$ws.Cells.Item(113, 11).Value = 3083.1667  # K113
$ws.Cells.Item(113, 13).Value = -913.1667000000002  # M113
$ws.Cells.Item(132, 8).Value = 2415.4  # H132
$ws.Cells.Item(132, 9).Value = 2749  # I132
$ws.Cells.Item(132, 11).Value = 8247  # K132
$ws.Cells.Item(132, 13).Value = -5717  # M132
$ws.Cells.Item(136, 8).Value = 5024.364  # H136
$ws.Cells.Item(136, 9).Value = 2500  # I136
$ws.Cells.Item(136, 11).Value = 7500  # K136
$ws.Cells.Item(136, 13).Value = -4950  # M136

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(48, 8).Value = 0  # H48
$ws.Cells.Item(48, 10).Value = 0  # J48
$ws.Cells.Item(48, 12).Value = 0  # L48
$ws.Cells.Item(48, 14).ClearContents()  # N48
$ws.Cells.Item(107, 8).Value = 542.5833  # H107
$ws.Cells.Item(107, 9).Value = 454.7143  # I107
$ws.Cells.Item(107, 11).Value = 1364.1429  # K107
$ws.Cells.Item(107, 13).Value = 555.8571000000002  # M107
$ws.Cells.Item(108, 8).Value = 55249.5  # H108
$ws.Cells.Item(108, 10).Value = 55249.5  # J108
$ws.Cells.Item(108, 12).Value = 55249.5  # L108
$ws.Cells.Item(108, 14).Value = -62929.5  # N108
$ws.Cells.Item(113, 8).Value = 1094.125  # H113
$ws.Cells.Item(113, 9).Value = 1030.6  # I113
$ws.Cells.Item(113, 11).Value = 3091.8  # K113
$ws.Cells.Item(113, 13).Value = -921.7999999999997  # M113
$ws.Cells.Item(123, 8).Value = 46153.453  # H123
$ws.Cells.Item(123, 10).Value = 46153.453  # J123
$ws.Cells.Item(123, 12).Value = 46153.453  # L123
$ws.Cells.Item(123, 14).Value = -55953.453  # N123
$ws.Cells.Item(132, 8).Value = 3474.2104  # H132
$ws.Cells.Item(132, 9).Value = 3363.7144  # I132
$ws.Cells.Item(132, 11).Value = 10091.1432  # K132
$ws.Cells.Item(132, 13).Value = -7561.143199999999  # M132
$ws.Cells.Item(136, 8).Value = 4145.684  # H136
$ws.Cells.Item(136, 9).Value = 4252  # I136
$ws.Cells.Item(136, 11).Value = 12756  # K136
$ws.Cells.Item(136, 13).Value = -10206  # M136
